# Voltage limit can now be edited
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 41: change the GB column (F) value from the hard-coded "23.4" placeholder
# to the editable "<value>" placeholder.
$ws.Range("F41").Value = "<value>"

# Row 42: new entry for the (now editable) voltage limit value.
$ws.Range("B42").Value = "SingleUseId69"
$ws.Range("C42").Value = "Tiny"
$ws.Range("D42").Value = "Right"
$ws.Range("E42").Value = "LTR"
$ws.Range("F42").Value = "1232<value>"
